$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.264129366914366, 1.448131862155378, 7.934253886187869, 2.816780766440276, 2.870363925838609, 22)
    3 = @(0.5477676633828543, 1.907452784230856, 14.64723191001851, 3.827170222242344, 3.881306689997472, 21)
    4 = @(-0.8445613698142435, 1.368339376243751, 7.290098563283141, 2.700018252398146, 2.631152370290416, 20)
    5 = @(-0.02661223054417471, 0.470512716437627, 0.5301842917739205, 0.7281375500370245, 0.7475904084393145, 19)
    6 = @(0.07547341295830622, 0.7121491052738879, 1.142006432124239, 1.068647010066579, 1.096882843470316, 18)
    7 = @(-0.2117013632410417, 0.6042393667413337, 0.8609082534243339, 0.9278514177519663, 0.9311801438848216, 17)
    8 = @(-0.06838949262952414, 0.4031547194043731, 0.4042733775863497, 0.6358249582914701, 0.6543231957670564, 15)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
